{"js": "// Update the worksheet's division-problem answers in place.\n// The document contains a single table; only every 4th row (0, 4, 8, 12, 16)\n// holds the five \"a\u00f7b=c, r\" answer cells per row, the rows in between are\n// blank spacer rows. We replace each answer cell's text with its new value,\n// addressing cells positionally (table-row-index, column-index) rather than\n// searching by old text, since several new values collide with other old\n// values elsewhere in the table.\n\nconst table = context.document.body.tables.items[0];\ntable.rows.load(\"items\");\nawait context.sync();\n\n// New values, keyed by the table's actual row index (including the blank\n// spacer rows) and column index (0-based), in document order.\nconst newValuesByRow = {\n  0: [\"80\u00f72=40, 0\", \"28\u00f76=4, 4\", \"48\u00f76=8, 0\", \"57\u00f72=28, 1\", \"16\u00f77=2, 2\"],\n  4: [\"55\u00f75=11, 0\", \"33\u00f74=8, 1\", \"91\u00f73=30, 1\", \"98\u00f72=49, 0\", \"64\u00f73=21, 1\"],\n  8: [\"16\u00f77=2, 2\", \"20\u00f74=5, 0\", \"96\u00f73=32, 0\", \"71\u00f73=23, 2\", \"89\u00f73=29, 2\"],\n  12: [\"48\u00f76=8, 0\", \"51\u00f74=12, 3\", \"23\u00f76=3, 5\", \"63\u00f77=9, 0\", \"45\u00f79=5, 0\"],\n  16: [\"12\u00f74=3, 0\", \"24\u00f73=8, 0\", \"36\u00f72=18, 0\", \"70\u00f75=14, 0\", \"85\u00f72=42, 1\"],\n};\n\nfor (const rowIndexStr of Object.keys(newValuesByRow)) {\n  const rowIndex = Number(rowIndexStr);\n  const values = newValuesByRow[rowIndex];\n  for (let col = 0; col < values.length; col++) {\n    const cell = table.getCell(rowIndex, col);\n    cell.value = values[col];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet's division-problem answers in place.\n# The document contains a single table; only every 4th row (1, 5, 9, 13, 17\n# in Word's 1-based Cell() indexing) holds the five \"a\u00f7b=c, r\" answer cells\n# per row, the rows in between are blank spacer rows. We replace each answer\n# cell's text with its new value, addressing cells positionally (table row,\n# column) rather than searching by old text, since several new values\n# collide with other old values elsewhere in the table.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newValuesByRow = @{\n    1  = @(\"80\u00f72=40, 0\", \"28\u00f76=4, 4\", \"48\u00f76=8, 0\", \"57\u00f72=28, 1\", \"16\u00f77=2, 2\")\n    5  = @(\"55\u00f75=11, 0\", \"33\u00f74=8, 1\", \"91\u00f73=30, 1\", \"98\u00f72=49, 0\", \"64\u00f73=21, 1\")\n    9  = @(\"16\u00f77=2, 2\", \"20\u00f74=5, 0\", \"96\u00f73=32, 0\", \"71\u00f73=23, 2\", \"89\u00f73=29, 2\")\n    13 = @(\"48\u00f76=8, 0\", \"51\u00f74=12, 3\", \"23\u00f76=3, 5\", \"63\u00f77=9, 0\", \"45\u00f79=5, 0\")\n    17 = @(\"12\u00f74=3, 0\", \"24\u00f73=8, 0\", \"36\u00f72=18, 0\", \"70\u00f75=14, 0\", \"85\u00f72=42, 1\")\n}\n\nforeach ($row in $newValuesByRow.Keys) {\n    $values = $newValuesByRow[$row]\n    for ($col = 1; $col -le $values.Count; $col++) {\n        $t.Cell($row, $col).Range.Text = $values[$col - 1]\n    }\n}\n"}
